$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Replace the tiny placeholder picture (the 1x1px "Submission Flowchart"
# image) with a plain hyperlink that points at the image's real URL on
# ura.gov.sg, displayed as visible link text.
# ---------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$rng = $shp.Range
$shp.Delete()
$hl = $d.Hyperlinks.Add($rng, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg", $null, $null, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg")

Write-Output "done"
